# Auto-generated edit script applying scheduled market-price updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 839
$ws.Range("I12").Value = 131.66667
$ws.Range("J12").Value = 1900
$ws.Range("K12").Value = 131.66667
$ws.Range("L12").Value = 1900
$ws.Range("M12").Value = 38.33332999999999
$ws.Range("N12").Value = -2240
$ws.Range("H70").Value = 21706660
$ws.Range("I70").Value = 14712567
$ws.Range("K70").Value = 44137701
$ws.Range("M70").Value = -44137431
$ws.Range("H73").Value = 21706660
$ws.Range("I73").Value = 14712567
$ws.Range("K73").Value = 44137701
$ws.Range("M73").Value = -44136765
$ws.Range("H112").Value = 6138.421
$ws.Range("J112").Value = 7981.4287
$ws.Range("L112").Value = 23944.2861
$ws.Range("N112").Value = -26160.2861
$ws.Range("H129").Value = 1315.6316
$ws.Range("I129").Value = 688.1
$ws.Range("K129").Value = 2064.3
$ws.Range("M129").Value = 2935.7
$ws.Range("H132").Value = 1893
$ws.Range("I132").Value = 1893
$ws.Range("K132").Value = 5679
$ws.Range("M132").Value = -3149
$ws.Range("H137").Value = 7420.222
$ws.Range("I137").Value = 5852.857
$ws.Range("K137").Value = 17558.571
$ws.Range("M137").Value = -15008.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2455432.8
$ws.Range("I32").Value = 2455432.8
$ws.Range("K32").Value = 2455432.8
$ws.Range("M32").Value = -2455145.8
$ws.Range("H61").Value = 50007540
$ws.Range("I61").Value = 1451.5
$ws.Range("J61").Value = 83344936
$ws.Range("K61").Value = 1451.5
$ws.Range("L61").Value = 83344936
$ws.Range("M61").Value = -1239.5
$ws.Range("N61").Value = -83345360
$ws.Range("H122").Value = 12303.519
$ws.Range("I122").Value = 18034.5
$ws.Range("J122").Value = 6131.6924
$ws.Range("K122").Value = 54103.5
$ws.Range("L122").Value = 18395.0772
$ws.Range("M122").Value = -51653.5
$ws.Range("N122").Value = -23295.0772
$ws.Range("H132").Value = 7382.421
$ws.Range("I132").Value = 4564.6665
$ws.Range("K132").Value = 13693.9995
$ws.Range("M132").Value = -11163.9995
$ws.Range("H136").Value = 50007540
$ws.Range("I136").Value = 1451.5
$ws.Range("J136").Value = 83344936
$ws.Range("K136").Value = 4354.5
$ws.Range("L136").Value = 250034808
$ws.Range("M136").Value = -1804.5
$ws.Range("N136").Value = -250039908

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10420225
$ws.Range("I20").Value = 18521168
$ws.Range("K20").Value = 18521168
$ws.Range("M20").Value = -18520921
$ws.Range("H86").Value = 42391.68
$ws.Range("I86").Value = 64450.25
$ws.Range("K86").Value = 64450.25
$ws.Range("M86").Value = -63327.25
$ws.Range("H89").Value = 42391.68
$ws.Range("I89").Value = 64450.25
$ws.Range("K89").Value = 322251.25
$ws.Range("M89").Value = -316635.25
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 89999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 89999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 89999
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -100279

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4459.381
$ws.Range("I16").Value = 1099.2858
$ws.Range("K16").Value = 1099.2858
$ws.Range("M16").Value = -812.2858000000001
$ws.Range("H32").Value = 1920
$ws.Range("J32").Value = 2250
$ws.Range("L32").Value = 2250
$ws.Range("N32").Value = -2882
$ws.Range("H35").Value = 482.33334
$ws.Range("I35").Value = 438.8
$ws.Range("J35").Value = 700
$ws.Range("K35").Value = 438.8
$ws.Range("L35").Value = 700
$ws.Range("M35").Value = -144.8
$ws.Range("N35").Value = -1288
$ws.Range("H42").Value = 25514
$ws.Range("J42").Value = 44500
$ws.Range("L42").Value = 44500
$ws.Range("N42").Value = -45686
$ws.Range("H58").Value = 9411.046
$ws.Range("I58").Value = 1614.2
$ws.Range("K58").Value = 1614.2
$ws.Range("M58").Value = -1411.2
$ws.Range("H86").Value = 12503080
$ws.Range("I86").Value = 20835666
$ws.Range("J86").Value = 4200
$ws.Range("K86").Value = 20835666
$ws.Range("L86").Value = 4200
$ws.Range("M86").Value = -20834543
$ws.Range("N86").Value = -6446
$ws.Range("H89").Value = 12503080
$ws.Range("I89").Value = 20835666
$ws.Range("J89").Value = 4200
$ws.Range("K89").Value = 104178330
$ws.Range("L89").Value = 21000
$ws.Range("M89").Value = -104172714
$ws.Range("N89").Value = -32232
$ws.Range("H99").Value = 11437.5
$ws.Range("I99").Value = 14299.2
$ws.Range("K99").Value = 14299.2
$ws.Range("M99").Value = -12801.2
$ws.Range("H110").Value = 92000
$ws.Range("J110").Value = 92000
$ws.Range("L110").Value = 92000
$ws.Range("N110").Value = -100180
$ws.Range("H113").Value = 4459.381
$ws.Range("I113").Value = 1099.2858
$ws.Range("K113").Value = 1099.2858
$ws.Range("M113").Value = 1070.7142
$ws.Range("H126").Value = 11437.5
$ws.Range("I126").Value = 14299.2
$ws.Range("K126").Value = 42897.60000000001
$ws.Range("M126").Value = -40427.60000000001
$ws.Range("H132").Value = 6747.1333
$ws.Range("I132").Value = 5275
$ws.Range("K132").Value = 15825
$ws.Range("M132").Value = -13295
$ws.Range("H136").Value = 9411.046
$ws.Range("I136").Value = 1614.2
$ws.Range("K136").Value = 4842.6
$ws.Range("M136").Value = -2292.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 830.7143
$ws.Range("I18").Value = 363.2
$ws.Range("K18").Value = 1089.6
$ws.Range("M18").Value = -920.5999999999999
$ws.Range("H86").Value = 2003
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2003
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 186857.22
$ws.Range("I113").Value = 558492.8
$ws.Range("K113").Value = 558492.8
$ws.Range("M113").Value = -556322.8
$ws.Range("H122").Value = 5177872
$ws.Range("I122").Value = 7246747.5
$ws.Range("K122").Value = 21740242.5
$ws.Range("M122").Value = -21737792.5
$ws.Range("H132").Value = 8827.916999999999
$ws.Range("I132").Value = 3302.2
$ws.Range("K132").Value = 9906.599999999999
$ws.Range("M132").Value = -7376.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4271.2856
$ws.Range("I7").Value = 2666.8
$ws.Range("K7").Value = 2666.8
$ws.Range("M7").Value = -2554.8
$ws.Range("H16").Value = 1500
$ws.Range("I16").Value = 1500
$ws.Range("K16").Value = 1500
$ws.Range("M16").Value = -1330
$ws.Range("H68").Value = 4891.25
$ws.Range("I68").Value = 5188.3335
$ws.Range("J68").Value = 4000
$ws.Range("K68").Value = 5188.3335
$ws.Range("L68").Value = 4000
$ws.Range("M68").Value = -4439.3335
$ws.Range("N68").Value = -5498
$ws.Range("H71").Value = 4891.25
$ws.Range("I71").Value = 5188.3335
$ws.Range("J71").Value = 4000
$ws.Range("K71").Value = 25941.6675
$ws.Range("L71").Value = 20000
$ws.Range("M71").Value = -22197.6675
$ws.Range("N71").Value = -27488
$ws.Range("H82").Value = 503138.16
$ws.Range("I82").Value = 1113908.5
$ws.Range("J82").Value = 3416.9092
$ws.Range("K82").Value = 1113908.5
$ws.Range("L82").Value = 3416.9092
$ws.Range("M82").Value = -1113547.5
$ws.Range("N82").Value = -4138.9092
$ws.Range("H85").Value = 503138.16
$ws.Range("I85").Value = 1113908.5
$ws.Range("J85").Value = 3416.9092
$ws.Range("K85").Value = 1113908.5
$ws.Range("L85").Value = 3416.9092
$ws.Range("M85").Value = -1112660.5
$ws.Range("N85").Value = -5912.9092
$ws.Range("H93").Value = 2682.111
$ws.Range("J93").Value = 2377.6667
$ws.Range("L93").Value = 2377.6667
$ws.Range("N93").Value = -4873.6667
$ws.Range("H100").Value = 4762.1904
$ws.Range("I100").Value = 2318.125
$ws.Range("K100").Value = 2318.125
$ws.Range("M100").Value = -1777.125
$ws.Range("H126").Value = 4271.2856
$ws.Range("I126").Value = 2666.8
$ws.Range("K126").Value = 8000.400000000001
$ws.Range("M126").Value = -5530.400000000001
$ws.Range("H132").Value = 13901289
$ws.Range("I132").Value = 35717548
$ws.Range("K132").Value = 107152644
$ws.Range("M132").Value = -107150114
$ws.Range("H136").Value = 9158.085999999999
$ws.Range("I136").Value = 3289.5715
$ws.Range("K136").Value = 9868.7145
$ws.Range("M136").Value = -7318.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17198.334
$ws.Range("J41").Value = 17198.334
$ws.Range("L41").Value = 17198.334
$ws.Range("N41").Value = -17978.334
$ws.Range("H122").Value = 14005390
$ws.Range("I122").Value = 17384774
$ws.Range("K122").Value = 52154322
$ws.Range("M122").Value = -52151872
$ws.Range("H132").Value = 27783420
$ws.Range("I132").Value = 40006150
$ws.Range("J132").Value = 4485.1816
$ws.Range("K132").Value = 120018450
$ws.Range("L132").Value = 13455.5448
$ws.Range("M132").Value = -120015920
$ws.Range("N132").Value = -18515.5448

Write-Host "Done. Sets:" 235 "Clears:" 4
